$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("dataSheet")

$ws.Range("A14").Value = "othervar"
$ws.Range("B14").Value = "something else"

$ws.Range("B14").Select()
